$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 6, pushing existing rows 6-15 down to 7-16
$ws.Rows("6:6").Insert()

# Populate the newly inserted row 6 with the new weekly record
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C6").Value = "Arica y Parinacota"
$ws.Range("D6").Value = 44547
$ws.Range("E6").Value = 15
$ws.Range("F6").Value = 100112003
$ws.Range("G6").Value = "Ajo"
$ws.Range("H6").Value = "Chino"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 300
$ws.Range("K6").Value = 19000
$ws.Range("L6").Value = 20000
$ws.Range("M6").Value = 19500
$ws.Range("N6").Value = "$/caja 10 kilos"
$ws.Range("O6").Value = "China"
$ws.Range("P6").Value = 1950
$ws.Range("Q6").Value = 10
$ws.Range("R6").Value = "Hortaliza"
